$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20; existing rows 20..58 shift down to 21..59.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record's data.
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44614
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = 100112012
$ws.Range("G20").Value = "Espinaca"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 1750
$ws.Range("N20").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 583
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = "Hortaliza"
